$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: mark "Wt" object as done (added sex column to B6)
$ws.Range("B6").Value = "y"

# Row 29 (M2): new "Questions" column entries
$ws.Range("E29").Value = "wt"
$ws.Range("F29").Value = "added sex column"

# Row 30 (UobsWtAge_hat): new "Questions" column entries
$ws.Range("E30").Value = "M1"
$ws.Range("F30").Value = "added sex column"

# Row 31 (mn_UobsWtAge_hat): new note
$ws.Range("E31").Value = "Combined fsh_control and srv_control to Fleet_control"

# Row 32 (suit_main): new note
$ws.Range("E32").Value = "Combined empirical selectivity"

# New row 33: LbyAge object + note
$ws.Range("A33").Value = "LbyAge"
$ws.Range("E33").Value = "Combined fsh_comp and srv_comp"

# New row 34: ConsumAge object
$ws.Range("A34").Value = "ConsumAge"

# Update the active selection to match where the author left off editing
$ws.Range("E33").Select()
